# Insert a new reference row (row 2) pointing to the POMS citation used for
# retroactive application, mirroring the existing "reference" note already
# in A1. This adds a new shared string + a new populated row without
# disturbing the position of the existing table rows (row 2 was
# previously blank, so nothing else needs to move).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 1)
$cell.Value = "https://secure.ssa.gov/poms.nsf/lnx/0200204030"

# Match the look of the other wrapped/top-aligned notes cells.
$cell.Font.Bold = $false
$cell.WrapText = $true
$cell.VerticalAlignment = -4160   # xlTop

# The row needs extra height to show the wrapped text (matches ht="34" in
# the saved file).
$ws.Rows.Item(2).RowHeight = 34

# Author's cursor ended up on the new cell after adding it.
[void]$ws.Range("A2").Select()
